# Pharmacy portal changes: add a "Batch numbe" column to Sheet1 between
# "medicine_type" (D) and "expiry_date" (old E, now F), with batch
# numbers "1b_10" / "2b_10" for the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at E, shifting expiry_date/dosage/qty/... right.
$ws.Columns("E:E").Insert()

# Populate the new column.
$ws.Range("E1").Value = "Batch numbe"
$ws.Range("E2").Value = "1b_10"
$ws.Range("E3").Value = "2b_10"

# Leave the same selection state the source workbook ended up with.
$ws.Range("E3").Select() | Out-Null
